$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.742069482803345
$ws.Range("B1").Value = 1.713589191436768
$ws.Range("C1").Value = 7.653044700622559
$ws.Range("D1").Value = 1.107338070869446
$ws.Range("E1").Value = 0.4277065992355347
